# The deck currently carries the "Integral" theme (dk2/lt2/accent1-6/hlink/
# folHlink recolored) on its one-and-only Slide Master/Theme pair. The
# target commit swaps that for the stock "Office Theme" palette (the same
# palette that PowerPoint ships by default), while dk1/lt1 stay black/white
# either way.
#
# Helper: pack three 0-255 byte values into the little-endian 0xBBGGRR
# integer that PowerPoint's ColorFormat/ThemeColor .RGB property expects.
function Get-BGR($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$cs = $theme.ThemeColorScheme

# Office Theme color scheme (target state), in theme color-index order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
$cs.Item(1).RGB  = Get-BGR 0x00 0x00 0x00   # dk1      -> 000000
$cs.Item(2).RGB  = Get-BGR 0xFF 0xFF 0xFF   # lt1      -> FFFFFF
$cs.Item(3).RGB  = Get-BGR 0x44 0x54 0x6A   # dk2      -> 44546A
$cs.Item(4).RGB  = Get-BGR 0xE7 0xE6 0xE6   # lt2      -> E7E6E6
$cs.Item(5).RGB  = Get-BGR 0x5B 0x9B 0xD5   # accent1  -> 5B9BD5
$cs.Item(6).RGB  = Get-BGR 0xED 0x7D 0x31   # accent2  -> ED7D31
$cs.Item(7).RGB  = Get-BGR 0xA5 0xA5 0xA5   # accent3  -> A5A5A5
$cs.Item(8).RGB  = Get-BGR 0xFF 0xC0 0x00   # accent4  -> FFC000
$cs.Item(9).RGB  = Get-BGR 0x44 0x72 0xC4   # accent5  -> 4472C4
$cs.Item(10).RGB = Get-BGR 0x70 0xAD 0x47   # accent6  -> 70AD47
$cs.Item(11).RGB = Get-BGR 0x05 0x63 0xC1   # hlink    -> 0563C1
$cs.Item(12).RGB = Get-BGR 0x95 0x4F 0x72   # folHlink -> 954F72
